$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 최종점수 (column K) for rows 2-5
$ws.Range("K2").Value = 60.8
$ws.Range("K3").Value = 58.4
$ws.Range("K4").Value = 51.6
$ws.Range("K5").Value = 48.8

# Update MACRO_SCORE (column N) for rows 2-5
$ws.Range("N2").Value = 54.77309453746771
$ws.Range("N3").Value = 54.77309453746771
$ws.Range("N4").Value = 54.77309453746771
$ws.Range("N5").Value = 54.77309453746771
